$d = $word.ActiveDocument

# 1) Apply strikethrough formatting to the three milestone bullets that are now
#    considered complete/deprioritised (both paragraph mark and run get <w:strike/>).
$strikeTargets = @(
    "Implement basic forms of communication for implemented mechanics.",
    "Finalize the tuning of the variable constants which influence the behavioural physics of both the ball hitting mechanics as well as the environmental aspects.",
    "Finalize ability mechanics to work smoothly within their respective level environment."
)

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text.Trim()
    foreach ($target in $strikeTargets) {
        if ($text -eq $target) {
            $p.Range.Font.StrikeThrough = 1
        }
    }
}

# 2) Remove the three Notes paragraphs that are no longer relevant, shifting the
#    remaining Notes bullets up (the diff reads as a cascading text replacement
#    but is semantically just these three paragraph deletions).
$deleteTargets = @(
    "Need to get ready for play testing by tomorrow!",
    "Through ball issue where target ball can go through",
    "Increase wall bounciness so ball doesn’t stop on wall"
)

foreach ($target in $deleteTargets) {
    foreach ($p in $d.Paragraphs) {
        $text = $p.Range.Text.Trim()
        if ($text -eq $target) {
            $p.Range.Delete()
            break
        }
    }
}
